$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "JavaFile TestCase"

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "Method Name"
$ws.Cells.Item(1, 2).Value = "Inputs"
$ws.Cells.Item(1, 3).Value = "Expected Output"
$ws.Cells.Item(1, 4).Value = "Expected Status Code"

# Data rows
$ws.Cells.Item(2, 1).Value = "checkMinBalance"
$ws.Cells.Item(2, 2).Value = "[5001]"
$ws.Cells.Item(2, 3).Value = "You have sufficient balance amount"
$ws.Cells.Item(2, 4).Value = 200

$ws.Cells.Item(3, 1).Value = "checkMinBalance"
$ws.Cells.Item(3, 2).Value = "[5000]"
$ws.Cells.Item(3, 3).Value = "You have sufficient balance amount"
$ws.Cells.Item(3, 4).Value = 200

$ws.Cells.Item(4, 1).Value = "checkMinBalance"
$ws.Cells.Item(4, 2).Value = "[4999]"
$ws.Cells.Item(4, 3).Value = "Your account balance amount is lesser than minimum balance"
$ws.Cells.Item(4, 4).Value = 200

$ws.Cells.Item(5, 1).Value = "checkMinBalance"
$ws.Cells.Item(5, 2).Value = "[-1001]"
$ws.Cells.Item(5, 3).Value = "Invalid balance amount: Balance amount should be a positive number"
$ws.Cells.Item(5, 4).Value = 400

# Header row styling: left/center alignment - apply to A1 first then copy format across so only
# a single new style entry is produced.
$headerFirst = $ws.Cells.Item(1, 1)
$headerFirst.HorizontalAlignment = -4131
$headerFirst.VerticalAlignment = -4108
$headerFirst.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row height
$ws.Rows.Item(1).RowHeight = 26.5

# Column widths (engine rounds to nearest 1/6 character; these are the closest
# achievable inputs to the target stored widths of 59.90625 / 18.81640625)
$ws.Columns.Item(3).ColumnWidth = 59
$ws.Columns.Item(4).ColumnWidth = 18

# Selection of header row (mirrors the saved view state)
$null = $ws.Range("A1:XFD1").Select()

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1
